# Auto-generated edit script: refresh market-price derived columns (H-N)
# on 'Moogle Profits' leve-crafting sheets, per scheduled runner commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 6153.375
$ws.Range("J112").Value = 7854.6665
$ws.Range("L112").Value = 23563.9995
$ws.Range("N112").Value = -25779.9995

$ws.Range("H113").Value = 4742.8887
$ws.Range("I113").Value = 4836
$ws.Range("K113").Value = 4836
$ws.Range("M113").Value = -1582

$ws.Range("H141").Value = 2715.1936
$ws.Range("I141").Value = 1570.84
$ws.Range("J141").Value = 7483.3335
$ws.Range("K141").Value = 4712.52
$ws.Range("L141").Value = 22450.0005
$ws.Range("M141").Value = 467.4800000000005
$ws.Range("N141").Value = -32810.00049999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2954.2273
$ws.Range("I61").Value = 1533.3636
$ws.Range("K61").Value = 1533.3636
$ws.Range("M61").Value = -1321.3636

$ws.Range("H74").Value = 15629960
$ws.Range("I74").Value = 26317594
$ws.Range("J74").Value = 9571.538
$ws.Range("K74").Value = 26317594
$ws.Range("L74").Value = 9571.538
$ws.Range("M74").Value = -26316720
$ws.Range("N74").Value = -11319.538

$ws.Range("H77").Value = 15629960
$ws.Range("I77").Value = 26317594
$ws.Range("J77").Value = 9571.538
$ws.Range("K77").Value = 131587970
$ws.Range("L77").Value = 47857.69
$ws.Range("M77").Value = -131583602
$ws.Range("N77").Value = -56593.69

$ws.Range("H132").Value = 3393.5
$ws.Range("I132").Value = 2174.7827
$ws.Range("K132").Value = 6524.348100000001
$ws.Range("M132").Value = -3994.348100000001

$ws.Range("H136").Value = 2954.2273
$ws.Range("I136").Value = 1533.3636
$ws.Range("K136").Value = 4600.0908
$ws.Range("M136").Value = -2050.0908

$ws.Range("H139").Value = 75224
$ws.Range("J139").Value = 75224
$ws.Range("L139").Value = 75224
$ws.Range("N139").Value = -85504

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1228.2693
$ws.Range("I20").Value = 1226.3125
$ws.Range("J20").Value = 1231.4
$ws.Range("K20").Value = 1226.3125
$ws.Range("L20").Value = 1231.4
$ws.Range("M20").Value = -979.3125
$ws.Range("N20").Value = -1725.4

$ws.Range("H64").Value = 1150.4445
$ws.Range("I64").Value = 1226.3334
$ws.Range("J64").Value = 998.6667
$ws.Range("K64").Value = 1226.3334
$ws.Range("L64").Value = 998.6667
$ws.Range("M64").Value = -1001.3334
$ws.Range("N64").Value = -1448.6667

$ws.Range("H67").Value = 1150.4445
$ws.Range("I67").Value = 1226.3334
$ws.Range("J67").Value = 998.6667
$ws.Range("K67").Value = 1226.3334
$ws.Range("L67").Value = 998.6667
$ws.Range("M67").Value = -446.3334
$ws.Range("N67").Value = -2558.6667

$ws.Range("H74").Value = 30956.666
$ws.Range("I74").Value = 27090
$ws.Range("K74").Value = 27090
$ws.Range("M74").Value = -26154

$ws.Range("H77").Value = 30956.666
$ws.Range("I77").Value = 27090
$ws.Range("K77").Value = 81270
$ws.Range("M77").Value = -76590

$ws.Range("H86").Value = 1830.238
$ws.Range("I86").Value = 1860.9286
$ws.Range("J86").Value = 1768.8572
$ws.Range("K86").Value = 1860.9286
$ws.Range("L86").Value = 1768.8572
$ws.Range("M86").Value = -737.9286
$ws.Range("N86").Value = -4014.8572

$ws.Range("H89").Value = 1830.238
$ws.Range("I89").Value = 1860.9286
$ws.Range("J89").Value = 1768.8572
$ws.Range("K89").Value = 9304.643
$ws.Range("L89").Value = 8844.286
$ws.Range("M89").Value = -3688.643
$ws.Range("N89").Value = -20076.286

$ws.Range("H95").Value = 48924.332
$ws.Range("J95").Value = 48924.332
$ws.Range("L95").Value = 48924.332
$ws.Range("N95").Value = -54416.332

$ws.Range("H140").Value = 154997
$ws.Range("J140").Value = 154997
$ws.Range("L140").Value = 154997
$ws.Range("N140").Value = -165357

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 64995
$ws.Range("J124").Value = 64995
$ws.Range("L124").Value = 64995
$ws.Range("N124").Value = -69905

$ws.Range("H141").Value = 249193.72
$ws.Range("J141").Value = 249193.72
$ws.Range("L141").Value = 249193.72
$ws.Range("N141").Value = -259553.72

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 2000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 2000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 6000
$ws.Range("N9").Value = -6448
$ws.Range("M9").ClearContents()

$ws.Range("H64").Value = 7544.3335
$ws.Range("J64").Value = 3975
$ws.Range("L64").Value = 11925
$ws.Range("N64").Value = -12465

$ws.Range("H67").Value = 7544.3335
$ws.Range("J67").Value = 3975
$ws.Range("L67").Value = 11925
$ws.Range("N67").Value = -13797

$ws.Range("H80").Value = 4998.154
$ws.Range("I80").Value = 4999.3335
$ws.Range("J80").Value = 4997.8
$ws.Range("K80").Value = 14998.0005
$ws.Range("L80").Value = 14993.4
$ws.Range("M80").Value = -14062.0005
$ws.Range("N80").Value = -16865.4

$ws.Range("H83").Value = 4998.154
$ws.Range("I83").Value = 4999.3335
$ws.Range("J83").Value = 4997.8
$ws.Range("K83").Value = 44994.0015
$ws.Range("L83").Value = 44980.2
$ws.Range("M83").Value = -40314.0015
$ws.Range("N83").Value = -54340.2

$ws.Range("H92").Value = 420.42856
$ws.Range("I92").Value = 186
$ws.Range("J92").Value = 596.25
$ws.Range("K92").Value = 558
$ws.Range("L92").Value = 1788.75
$ws.Range("M92").Value = 690
$ws.Range("N92").Value = -4284.75

$ws.Range("H133").Value = 5026.815
$ws.Range("I133").Value = 4146.2856
$ws.Range("J133").Value = 5335
$ws.Range("K133").Value = 12438.8568
$ws.Range("L133").Value = 16005
$ws.Range("M133").Value = -7378.856800000001
$ws.Range("N133").Value = -26125

$ws.Range("H134").Value = 808.0769
$ws.Range("I134").Value = 808.0769
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2424.2307
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 2645.7693
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2728.8
$ws.Range("I122").Value = 1365.3334
$ws.Range("K122").Value = 4096.0002
$ws.Range("M122").Value = -1646.0002

$ws.Range("H129").Value = 32333.334
$ws.Range("J129").Value = 32333.334
$ws.Range("L129").Value = 32333.334
$ws.Range("N129").Value = -42333.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4142.24
$ws.Range("I16").Value = 3850.611
$ws.Range("J16").Value = 4892.143
$ws.Range("K16").Value = 3850.611
$ws.Range("L16").Value = 4892.143
$ws.Range("M16").Value = -3680.611
$ws.Range("N16").Value = -5232.143

$ws.Range("H93").Value = 988.6667
$ws.Range("I93").Value = 770.125
$ws.Range("K93").Value = 770.125
$ws.Range("M93").Value = 477.875

$ws.Range("H122").Value = 6310.3
$ws.Range("J122").Value = 8501.25
$ws.Range("L122").Value = 25503.75
$ws.Range("N122").Value = -30403.75

$ws.Range("H132").Value = 1980.7179
$ws.Range("I132").Value = 1235.1666
$ws.Range("K132").Value = 3705.4998
$ws.Range("M132").Value = -1175.4998

$ws.Range("H136").Value = 3025.2258
$ws.Range("I136").Value = 2289.5
$ws.Range("K136").Value = 6868.5
$ws.Range("M136").Value = -4318.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 588.9524
$ws.Range("I113").Value = 589
$ws.Range("J113").Value = 588.8333
$ws.Range("K113").Value = 1767
$ws.Range("L113").Value = 1766.4999
$ws.Range("M113").Value = 403
$ws.Range("N113").Value = -6106.4999

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 3014.7437
$ws.Range("I132").Value = 2130.743
$ws.Range("K132").Value = 6392.228999999999
$ws.Range("M132").Value = -3862.228999999999

Write-Output "Updated 39 rows across 8 sheets."